$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet and update title header text
$ws.Name = "Through 2022-02-09"
$ws.Range("B1").Value = "February 2022 (through February 09)"

# Apply cell value changes / additions for 2022-02-17 data
$ws.Range("B3").Value = 2
$ws.Range("H3").Value = 1
$ws.Range("L3").Value = 3

$ws.Range("D8").Value = 7

$ws.Range("D10").Value = 1
$ws.Range("L10").Value = 1

$ws.Range("D12").Value = 4

$ws.Range("D15").Value = 4

$ws.Range("D18").Value = 1

$ws.Range("F23").Value = 1

$ws.Range("D42").Value = 1

$ws.Range("D46").Value = 1

$ws.Range("D55").Value = 2

$ws.Range("D61").Value = 1

$ws.Range("F66").Value = 1

$ws.Range("F67").Value = 1

$ws.Range("D83").Value = 3

$ws.Range("H84").Value = 1

$ws.Range("J85").Value = 1
